$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format before writing numeric-looking strings,
# so values like "601.98" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '72.258.30'
$ws.Range("E2").Value = '  +1.60%  '
$ws.Range("D3").Value = '2.636.21'
$ws.Range("E3").Value = '  +0.94%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '601.98'
$ws.Range("E5").Value = '  -0.79%  '
$ws.Range("D6").Value = '180.03'
$ws.Range("E6").Value = '  -0.68%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("E8").Value = '  +0.46%  '
$ws.Range("E9").Value = '  +3.71%  '
$ws.Range("D10").Value = '2.634.72'
$ws.Range("E10").Value = '  +0.90%  '
$ws.Range("E11").Value = '  +1.31%  '
$ws.Range("D12").Value = '0.359'
$ws.Range("E12").Value = '  +2.71%  '
$ws.Range("E13").Value = '  -0.03%  '
$ws.Range("E14").Value = '  +2.80%  '
$ws.Range("D15").Value = '3.116.64'
$ws.Range("E15").Value = '  +1.50%  '
$ws.Range("D16").Value = '72.208.51'
$ws.Range("E16").Value = '  +1.60%  '
$ws.Range("D17").Value = '26.61'
$ws.Range("E17").Value = '  -1.01%  '
$ws.Range("D18").Value = '2.644.67'
$ws.Range("E18").Value = '  +0.44%  '
$ws.Range("D19").Value = '11.93'
$ws.Range("E19").Value = '  +4.14%  '
$ws.Range("D20").Value = '379.43'
$ws.Range("E20").Value = '  +0.25%  '
$ws.Range("E21").Value = '  -0.43%  '
$ws.Range("D22").Value = '4.19'
$ws.Range("E22").Value = '  -0.28%  '
$ws.Range("D23").Value = '2.09'
$ws.Range("E23").Value = '  +11.06%  '
$ws.Range("D24").Value = '73.25'
$ws.Range("E24").Value = '  +1.51%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").Value = '4.39'
$ws.Range("E26").Value = '  -0.66%  '
$ws.Range("D27").Value = '10.10'
$ws.Range("E27").Value = '  +3.82%  '
$ws.Range("D28").Value = '2.782.65'
$ws.Range("E28").Value = '  +1.28%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("D30").Value = '0.0₃0955'
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").Value = '524.23'
$ws.Range("E31").Value = '  -1.51%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '8.15'
$ws.Range("E32").Value = '  +0.47%  '
$ws.Range("E33").Value = '  -1.14%  '
$ws.Range("E34").Value = '  -0.82%  '
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("D36").Value = '164.98'
$ws.Range("E36").Value = '  +0.63%  '
$ws.Range("D37").Value = '19.34'
$ws.Range("E37").Value = '  +0.66%  '
$ws.Range("D38").Value = '0.113'
$ws.Range("E38").Value = '  -5.76%  '
$ws.Range("D39").Value = '19.08'
$ws.Range("E39").Value = '  +0.72%  '
$ws.Range("E40").Value = '  +1.43%  '
$ws.Range("D41").Value = '1.85'
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("E42").Value = '  +3.73%  '
$ws.Range("D43").Value = '5.07'
$ws.Range("E43").Value = '  -0.38%  '
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("E45").Value = '  +0.77%  '
$ws.Range("D46").Value = '39.33'
$ws.Range("E46").Value = '  -3.04%  '
$ws.Range("D47").Value = '151.08'
$ws.Range("E47").Value = '  -2.23%  '
$ws.Range("D48").Value = '3.72'
$ws.Range("E48").Value = '  +1.18%  '
$ws.Range("E49").Value = '  +2.01%  '
$ws.Range("D50").Value = '1.70'
$ws.Range("E50").Value = '  +1.34%  '
$ws.Range("D51").Value = '0.0₆0262'
$ws.Range("E51").Value = '  -3.01%  '

# Restore the original (default) style on column D so no stray number format
# metadata remains (matches the source workbook, which has no explicit style
# on these cells).
$ws.Range("D2:D51").Style = "Normal"
